$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the Cypher/Neo4j MATCH query text that drives this test case into A2
# (wrap-text style already present on A2 via the "Normal 2" style index).
$query = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Lung adenocarcinoma'] RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"

$ws.Range("A2").Value = $query

# The long wrapped query text makes the row taller once it is entered.
$ws.Rows("2:2").RowHeight = 87

# Selection moves to cover the new query cell, ending on A6.
$ws.Range("A2:A6").Select()
